$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the original D/J/K/L/M/P values for rows 2-31 before overwriting,
# since the edit re-distributes each row's data to a different row (a permutation).
$orig = @{}
$orig[2] = @(44223, 80, 2500, 3000, 2781, 927)
$orig[3] = @(44627, 78, 3500, 3500, 3500, 1167)
$orig[4] = @(44225, 56, 3000, 3000, 3000, 1000)
$orig[5] = @(44557, 104, 2000, 2500, 2260, 753)
$orig[6] = @(44340, 54, 3000, 3000, 3000, 1000)
$orig[7] = @(44537, 88, 2000, 2200, 2091, 697)
$orig[8] = @(44224, 67, 3000, 3000, 3000, 1000)
$orig[9] = @(44165, 68, 3000, 3000, 3000, 1000)
$orig[10] = @(44845, 80, 2500, 2500, 2500, 833)
$orig[11] = @(44935, 78, 3000, 3000, 3000, 1000)
$orig[12] = @(44222, 45, 3000, 3000, 3000, 1000)
$orig[13] = @(44804, 85, 3000, 3000, 3000, 1000)
$orig[14] = @(44193, 70, 3000, 3000, 3000, 1000)
$orig[15] = @(44574, 50, 3000, 3000, 3000, 1000)
$orig[16] = @(44389, 81, 2800, 3000, 2889, 963)
$orig[17] = @(44260, 60, 3500, 3500, 3500, 1167)
$orig[18] = @(44242, 95, 2500, 3000, 2737, 912)
$orig[19] = @(44187, 65, 3000, 3000, 3000, 1000)
$orig[20] = @(44166, 45, 2500, 2500, 2500, 833)
$orig[21] = @(44390, 50, 3000, 3000, 3000, 1000)
$orig[22] = @(44243, 45, 3000, 3000, 3000, 1000)
$orig[23] = @(44937, 68, 3500, 3500, 3500, 1167)
$orig[24] = @(44179, 78, 3000, 3000, 3000, 1000)
$orig[25] = @(44221, 50, 2500, 2500, 2500, 833)
$orig[26] = @(44669, 92, 2500, 3000, 2755, 918)
$orig[27] = @(44292, 40, 3000, 3000, 3000, 1000)
$orig[28] = @(44559, 68, 2000, 2000, 2000, 667)
$orig[29] = @(44536, 125, 2200, 2200, 2200, 733)
$orig[30] = @(44756, 104, 2800, 3000, 2904, 968)
$orig[31] = @(44291, 45, 3000, 3000, 3000, 1000)

# Column order in the snapshot arrays: D, J, K, L, M, P -> columns 4, 10, 11, 12, 13, 16
$targetCols = @(4, 10, 11, 12, 13, 16)

# Map: target row -> source row whose original values it should receive
$rowMap = @{}
$rowMap[2] = 7
$rowMap[3] = 19
$rowMap[4] = 2
$rowMap[5] = 4
$rowMap[6] = 17
$rowMap[7] = 5
$rowMap[8] = 16
$rowMap[9] = 13
$rowMap[10] = 21
$rowMap[11] = 11
$rowMap[12] = 8
$rowMap[13] = 9
$rowMap[14] = 24
$rowMap[15] = 31
$rowMap[16] = 27
$rowMap[17] = 29
$rowMap[18] = 30
$rowMap[19] = 12
$rowMap[20] = 26
$rowMap[21] = 20
$rowMap[22] = 23
$rowMap[23] = 22
$rowMap[24] = 6
$rowMap[25] = 18
$rowMap[26] = 25
$rowMap[27] = 28
$rowMap[28] = 10
$rowMap[29] = 3
$rowMap[30] = 15
$rowMap[31] = 14

foreach ($targetRow in ($rowMap.Keys | Sort-Object)) {
    $srcRow = $rowMap[$targetRow]
    $vals = $orig[$srcRow]
    for ($i = 0; $i -lt $targetCols.Length; $i++) {
        $ws.Cells.Item($targetRow, $targetCols[$i]).Value = $vals[$i]
    }
}
